$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows 50-132 down to 51-133.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with its data.
$ws.Range("A50").Value2 = 5
$ws.Range("B50").Value2 = "Macroferia Regional de Talca"
$ws.Range("C50").Value2 = "Maule"
$ws.Range("D50").Value2 = 44915
$ws.Range("E50").Value2 = 7
$ws.Range("F50").Value2 = 100112022
$ws.Range("G50").Value2 = "Arveja Verde"
$ws.Range("H50").Value2 = "Sin especificar"
$ws.Range("I50").Value2 = "Primera"
$ws.Range("J50").Value2 = 500
$ws.Range("K50").Value2 = 20000
$ws.Range("L50").Value2 = 20000
$ws.Range("M50").Value2 = 20000
$ws.Range("N50").Value2 = "`$/saco 25 kilos"
$ws.Range("O50").Value2 = "Carahue"
$ws.Range("P50").Value2 = 800
$ws.Range("Q50").Value2 = 25
$ws.Range("R50").Value2 = "Hortaliza"
